$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Overlap degree" cell updates (tool/vulnerability coverage matrix) ---
$ws.Range("B3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("D5").Value = 1

# --- View state: zoom to 150% and move the active selection to I3 ---
$win = $excel.ActiveWindow
$win.Zoom = 150
$ws.Range("I3").Select()
